$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.247.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.907.15'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.38%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5255'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3815'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07304'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("E10").Value = '  +1.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9058'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08173'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '96.51'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("E14").Value = '  +1.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.610.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -15.54%  '
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008689'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.77'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.090.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.128'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("E22").Value = '  +1.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.359'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.25'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '116.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.848'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.868'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09254'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8307'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05069'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.230'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.91%  '
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.355'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.738'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5772'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02008'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.081'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.132'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.594'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.39'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1523'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.38%  '
$ws.Range("E45").Value = '  +1.00%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.645'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '64.66'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06054'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.69%  '
